$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New changelog entry: row 10 gets a date + the new "Routing." note that
#    previously only existed as an empty placeholder row.
# ---------------------------------------------------------------------------

# Give A10 the same "date" formatting already used by A2/A3/A6 (numFmt 14,
# the light-blue fill + border), then stamp the serial date value.
$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 44793

# New note text in B10 (this becomes a brand-new shared string).
$ws.Range("B10").Value = "Routing. Almost finished routing for the VDP. Moved back the footprint for the 14mhz oscillator to THT as it is easier to find outside US/Europe and I have a bunch laying around here. Organized all components into the PCB. Removed the holed in the cartridge so I can better measure them in the standard MSX cartridge. Changed the RGB connector footprint to a shorter version so we can fit it in the standard size cartridge. Fixed multiple footprints to 0805 as we will be standardizing in that size. "

# Row 10 only needs 45pts now that column B is much wider (see below).
$ws.Rows.Item(10).RowHeight = 45

# ---------------------------------------------------------------------------
# 2. Column B becomes much wider (one long note per row) so every note cell
#    can wrap instead of needing a tall row. Widen the column, then give the
#    whole note column (header included) wrapped text, and shrink the two
#    rows that used to rely on tall non-wrapped rows.
# ---------------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 156

# B2:B19 already contains text; B2/B3 already wrap (copy that look onto the
# rest of the column so every note cell - old and new - shares one style).
$ws.Range("B2").Copy()
$ws.Range("B4:B19").PasteSpecial(-4122)

# The header cell B1 keeps its bold/fill/border look, just adds wrap.
$ws.Range("B1").WrapText = $true

# Rows 2 & 3 can shrink now that the text wraps in a much wider column.
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 75

# ---------------------------------------------------------------------------
# 3. Selection moves from the now-filled B10 to the next empty note cell.
# ---------------------------------------------------------------------------
[void]$ws.Range("B11").Select()
